$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.357.52"
$ws.Range("E2").Value = "  -3.80%  "
$ws.Range("D3").Value = "2.646.90"
$ws.Range("E3").Value = "  -1.57%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'521.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.66%  "
$ws.Range("D6").Value = "'144.89"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("E8").Value = "  -1.10%  "
$ws.Range("D9").Value = "'6.74"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.67%  "
$ws.Range("E10").Value = "  -2.97%  "
$ws.Range("D11").Value = "'0.339"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("E12").Value = "  +1.65%  "
$ws.Range("D13").Value = "3.114.96"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").Value = "58.364.18"
$ws.Range("E14").Value = "  -3.72%  "
$ws.Range("D15").Value = "'20.93"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("E16").Value = "  -1.44%  "
$ws.Range("D17").Value = "2.652.70"
$ws.Range("E17").Value = "  -1.85%  "
$ws.Range("D18").Value = "'338.98"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.10%  "
$ws.Range("E19").Value = "  -2.74%  "
$ws.Range("D20").Value = "'10.45"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "'6.34"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("E22").Value = "  +0.22%  "
$ws.Range("D23").Value = "'64.36"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.08%  "
$ws.Range("D24").Value = "'0.426"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +1.14%  "
$ws.Range("E25").Value = "  -2.08%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("D27").Value = "0.0₃0799"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("D28").Value = "'7.14"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.66%  "
$ws.Range("D29").Value = "'6.68"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -2.92%  "
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "'152.38"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "'18.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -1.96%  "
$ws.Range("D34").Value = "'4.16"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -5.66%  "
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("D37").Value = "'0.870"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.98%  "
$ws.Range("E38").Value = "  -0.42%  "
$ws.Range("E39").Value = "  -5.10%  "
$ws.Range("D40").Value = "'3.63"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +0.38%  "
$ws.Range("D42").Value = "'0.610"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'274.03"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.06%  "
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("D45").Value = "'19.45"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.37%  "
$ws.Range("B46").Value = "Hedera"
$ws.Range("C46").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D46").Value = "'0.0538"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.41%  "
$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").Value = "'10.63"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.54%  "
$ws.Range("D48").Value = "2.039.76"
$ws.Range("E48").Value = "  -4.79%  "
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "'4.69"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -4.43%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0229"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.48%  "
$ws.Range("D51").Value = "'18.30"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.54%  "
